# ---------------------------------------------------------------------------
# Applies the OOXML diff: isolates accented capital letters (O/I -> Ó/Í) that
# Word's AutoCorrect/spell-fix turned into their own runs, and wraps the
# {{codigo}} merge field in literal parentheses.
#
# Each edit below:
#   1. Locates the target text with Find.Execute (exact, case sensitive).
#   2. Mutates the text in place (replacing a letter, or inserting literal
#      characters).
#   3. "Stamps" a single character of formatting (toggle off/on, or on/off)
#      on each newly introduced boundary so the engine keeps it as its own
#      <w:r> run instead of silently re-merging it with its neighbour, while
#      leaving the final, saved formatting identical to the original.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-Boundary($range, [bool]$finalBold) {
    # Forces Word to keep $range (normally a single character) as an
    # independent run by toggling Bold away from, then back to, the value it
    # must end up with. The net formatting change is therefore zero.
    if ($finalBold) {
        $range.Bold = 0
        $range.Bold = 1
    } else {
        $range.Bold = 1
        $range.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# 1) "EVALUACION DE COMPATIBILIDAD DE USO" -> "EVALUACI" + "Ó" + "N DE ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("EVALUACION DE COMPATIBILIDAD DE USO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start
$accent = $d.Range($start + 8, $start + 9)
$accent.Text = "Ó"
Split-Boundary $accent $false

# ---------------------------------------------------------------------------
# 2) "DIRECCION                    " -> "DIRECCI" + "Ó" + "N                    "
#
# This run sits right before four other runs that carry byte-identical
# run-properties (same rFonts/sz/szCs/lang, only their w:rsidR differs), so
# mutating any text inside the first run causes the engine to coalesce all
# five of them into a single run. To land on the run layout the diff
# describes -- only run #1 split three ways, runs #2-#5 left untouched --
# every original run boundary has to be re-asserted afterwards, in addition
# to the two new ones introduced by the accented letter.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("DIRECCION                    ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start

$accent = $d.Range($start + 7, $start + 8)
$accent.Text = "Ó"

# (segStart, segEnd) offsets from $start, covering the original run #1
# split (DIRECCI / Ó / N + spaces) plus the original runs #2-#5 boundaries
# ("    ", "         ", " ", ":   ") which must stay independent runs.
$direccionSegments = @(
    @(0, 7),
    @(7, 8),
    @(8, 29),
    @(29, 33),
    @(33, 42),
    @(42, 43),
    @(43, 47)
)
foreach ($seg in $direccionSegments) {
    $segRange = $d.Range($start + $seg[0], $start + $seg[1])
    Split-Boundary $segRange $false
}

# ---------------------------------------------------------------------------
# 3) "ACTIVIDADES ... EN EL INDICE DE USOS ..." -> "... EN EL " + "Í" + "NDICE ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("ACTIVIDADES QUE REGISTRAN EN EL INDICE DE USOS DE ACTIVIDADES URBANAS CONFORME LO ESTABLECE LA ORD. 933-MML:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start
$accent = $d.Range($start + 32, $start + 33)
$accent.Text = "Í"
Split-Boundary $accent $false

# ---------------------------------------------------------------------------
# 4) "{{actividad}}{{codigo}}" -> "{{actividad}}" + "(" + "{{codigo}}" + ")"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("{{actividad}}{{codigo}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start
$end = $r.End

$openAt = $start + 13
$openIns = $d.Range($openAt, $openAt)
$openIns.InsertBefore("(")

$closeAt = $end + 1
$closeIns = $d.Range($closeAt, $closeAt)
$closeIns.InsertBefore(")")

$openRange = $d.Range($openAt, $openAt + 1)
Split-Boundary $openRange $true
$closeRange = $d.Range($closeAt, $closeAt + 1)
Split-Boundary $closeRange $true

# ---------------------------------------------------------------------------
# 5) "ZONIFICACION" -> "ZONIFICACI" + "Ó" + "N"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("ZONIFICACION", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start
$accent = $d.Range($start + 10, $start + 11)
$accent.Text = "Ó"
Split-Boundary $accent $true

# ---------------------------------------------------------------------------
# 6) "CODIGO" -> "C" + "Ó" + "DIGO"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("CODIGO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start
$accent = $d.Range($start + 1, $start + 2)
$accent.Text = "Ó"
Split-Boundary $accent $false

Write-Host "All edits applied."
